# Update the HCP path values on the test_case_Evan sheet to reflect the new
# (parallelized) CSV generation output, and move the active selection to
# where the user left off reviewing the results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the "PATH_HCP2" column (I) should now point at the
# NDARINV7Y7JEGPW scan instead of the NDARINV02EBX0JJ one.
$ws.Range("I2").Value = "data/hcp_comm_det_damien/cub-sub-NDARINV7Y7JEGPW_FNL_preproc_v2_Atlas_SMOOTHED_1.7.dtseries.nii_10_minutes_of_data_at_FD_0.2.dconn.nii_to_Merged_HCP_best80_dtseries.conc_AVG.dconn.dscalar.nii"

# Row 3: both the "PATH_HCP" (H) and "PATH_HCP2" (I) columns move off the
# NDARINV0U23M45G scan - H3 now matches H2's scan (NDARINV02EBX0JJ) and I3
# matches the same NDARINV7Y7JEGPW scan used above for I2.
$ws.Range("H3").Value = "data/hcp_comm_det_damien/cub-sub-NDARINV02EBX0JJ_FNL_preproc_v2_Atlas_SMOOTHED_1.7.dtseries.nii_10_minutes_of_data_at_FD_0.2.dconn.nii_to_Merged_HCP_best80_dtseries.conc_AVG.dconn.dscalar.nii"
$ws.Range("I3").Value = "data/hcp_comm_det_damien/cub-sub-NDARINV7Y7JEGPW_FNL_preproc_v2_Atlas_SMOOTHED_1.7.dtseries.nii_10_minutes_of_data_at_FD_0.2.dconn.nii_to_Merged_HCP_best80_dtseries.conc_AVG.dconn.dscalar.nii"

# Leave the selection where the reviewer ended up checking the test output.
$ws.Range("H14").Select()
